# A new weekly price record was inserted immediately before the existing
# row 286 (2023-01-25 / 25-01-2023, "Segunda" quality, volume 200, prices
# 17000 / $1214 per Kg). All the rows that used to be 286..333 shift down
# by one (now 287..334); row 334 therefore now holds the data that used to
# live in row 333.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 286..333 down to 287..334, leaving a blank row 286 behind
# (Excel copies formatting from the row above, same as an interactive
# right-click > Insert).
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(286, 1).Value  = 5
$ws.Cells.Item(286, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(286, 3).Value  = "Maule"
$ws.Cells.Item(286, 4).Value  = 44951
$ws.Cells.Item(286, 5).Value  = 7
$ws.Cells.Item(286, 6).Value  = "Fruta"
$ws.Cells.Item(286, 7).Value  = 100108
$ws.Cells.Item(286, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(286, 9).Value  = 100108005
$ws.Cells.Item(286, 10).Value = "Piña"
$ws.Cells.Item(286, 11).Value = "Caramelo"
$ws.Cells.Item(286, 12).Value = "Segunda"
$ws.Cells.Item(286, 13).Value = 200
$ws.Cells.Item(286, 14).Value = 17000
$ws.Cells.Item(286, 15).Value = 17000
$ws.Cells.Item(286, 16).Value = 17000
$ws.Cells.Item(286, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(286, 18).Value = "Ecuador"
$ws.Cells.Item(286, 19).Value = 1214
$ws.Cells.Item(286, 20).Value = 14
